$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 42

    # Trade #
    $ws.Cells.Item($row, 1).Value2 = 41

    # Date - leading apostrophe forces text entry so Excel doesn't
    # reinterpret the ISO-ish string as a date serial number.
    $cB = $ws.Cells.Item($row, 2)
    $cB.Value2 = "'2026-02-16"
    $cB.Style = "Normal"

    # Time (kept as plain text - not date/time-like enough to get
    # auto-converted by the engine).
    $ws.Cells.Item($row, 3).Value2 = "22:56:12"

    # Strategy
    $ws.Cells.Item($row, 4).Value2 = "base_strategy"

    # Side
    $ws.Cells.Item($row, 5).Value2 = "DOWN"

    # Entry Price
    $ws.Cells.Item($row, 6).Value2 = 49.999998

    # Exit Price - blank in source row; use the quote-prefix trick to
    # materialize an empty-text cell (same shared empty string as the
    # rest of the column) instead of leaving the cell absent.
    $cG = $ws.Cells.Item($row, 7)
    $cG.Value2 = "'"
    $cG.Style = "Normal"

    # Status
    $ws.Cells.Item($row, 8).Value2 = "OPEN"

    # P&L %
    $ws.Cells.Item($row, 9).Value2 = 0

    # P&L $
    $ws.Cells.Item($row, 10).Value2 = 0

    # Capital After
    $ws.Cells.Item($row, 11).Value2 = 100

    # Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value2 = 0

    # Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value2 = 0

    # Confidence
    $ws.Cells.Item($row, 14).Value2 = 0.6

    # Entry Reason
    $ws.Cells.Item($row, 15).Value2 = "Normal spread capture: 19600 bps"

    # Exit Reason - blank, same treatment as Exit Price above.
    $cP = $ws.Cells.Item($row, 16)
    $cP.Value2 = "'"
    $cP.Style = "Normal"

    # Duration (min)
    $ws.Cells.Item($row, 17).Value2 = 0
}
